$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the image file names from .jpg to .BMP (rows 2-4, column A)
$ws.Range("A2").Value = "images/Face1.BMP"
$ws.Range("A3").Value = "images/Face2.BMP"
$ws.Range("A4").Value = "images/Face3.BMP"

# Adjust column widths to match the new content layout
$ws.Columns.Item(1).ColumnWidth = 22.166666666666668
$ws.Columns.Item(2).ColumnWidth = 12.833333333333334
$ws.Columns.Item(3).ColumnWidth = 12.333333333333334

# Update the active selection
$ws.Range("B5").Select()
